$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): insert new "Responsable" column, shift Nombre/Fecha over ---
$ws.Range("B1").Value = "Responsable"
$ws.Range("C1").Value = "Nombre"
$ws.Range("D1").Value = "Fecha"

# --- Row 2: new data row ---
$ws.Range("A2").Value = "Registro"
$ws.Range("B2").Value = "Asesor"
$ws.Range("C2").Value = "'32452"
$ws.Range("D2").Value = "'2025-02-06"

# --- Row 3: new data row ---
$ws.Range("A3").Value = "Pedido"
$ws.Range("B3").Value = "Asesor"
$ws.Range("C3").Value = "asdsad"
$ws.Range("D3").Value = "'2025-02-06"
